# Update the "修改时间" (modified time) stamps across all three portfolio
# sheets, bumping the sync timestamp from 202509211440 to 202509211451
# (Web UI sync recorded at 2025-09-21 06:51).
#
# The timestamp is stored as text, so a leading apostrophe is used to keep
# Excel from reinterpreting the digit string as a number.

$wb = $excel.ActiveWorkbook

$newTimestamp = "'202509211451"

# Sheet "大智投资组合": timestamp lives in column E, data rows 2-9
$ws1 = $wb.Worksheets.Item("大智投资组合")
$ws1.Range("E2:E9").Value = $newTimestamp

# Sheet "大成投资组合": timestamp lives in column E, data rows 2-11
$ws2 = $wb.Worksheets.Item("大成投资组合")
$ws2.Range("E2:E11").Value = $newTimestamp

# Sheet "我的投资组合": timestamp lives in column G, data rows 2-13
$ws3 = $wb.Worksheets.Item("我的投资组合")
$ws3.Range("G2:G13").Value = $newTimestamp
